# Adding creator read access for new complaint
# Mirrors existing "Case File - creator read access" rule (row 26) but for
# the COMPLAINT object type, appended as a new row right after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of the last existing rule row (row 26, columns A-H)
# into the new row 27 so the new row matches the table's look & feel.
$ws.Range("A26:H26").Copy()
$ws.Range("A27").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(27).RowHeight = $ws.Rows.Item(26).RowHeight

# Fill in the values for the new "Complaint - creator read access" rule.
$ws.Range("B27").Value = "Complaint - creator read access"
$ws.Range("C27").Value = "COMPLAINT"
$ws.Range("H27").Value = "reader, creator"
